# DSM Scheduled Flights vs actual.xlsx
# Append 9 new daily rows (2022-01-15 .. 2022-01-23) below the existing
# data block (which ends at row 649), following the same pattern as the
# prior rows: DateTime | Scheduled flights | Tracked flights | =C/B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style/template row to clone formatting from (last existing data row).
$templateRow = 649

# New rows: DateTime, Scheduled flights (B), Tracked flights (C)
$newRows = @(
    @{ Row = 650; Date = "2022-01-15"; Scheduled = 50; Tracked = 41 },
    @{ Row = 651; Date = "2022-01-16"; Scheduled = 53; Tracked = 48 },
    @{ Row = 652; Date = "2022-01-17"; Scheduled = 60; Tracked = 54 },
    @{ Row = 653; Date = "2022-01-18"; Scheduled = 55; Tracked = 54 },
    @{ Row = 654; Date = "2022-01-19"; Scheduled = 55; Tracked = 51 },
    @{ Row = 655; Date = "2022-01-20"; Scheduled = 75; Tracked = 71 },
    @{ Row = 656; Date = "2022-01-21"; Scheduled = 70; Tracked = 66 },
    @{ Row = 657; Date = "2022-01-22"; Scheduled = 48; Tracked = 43 },
    @{ Row = 658; Date = "2022-01-23"; Scheduled = 57; Tracked = 54 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: DateTime (stored as text, matching the existing column).
    $ws.Cells.Item($templateRow, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = $r.Date

    # Column B: Scheduled flights.
    $ws.Cells.Item($templateRow, 2).Copy($ws.Cells.Item($row, 2))
    $ws.Cells.Item($row, 2).Value = $r.Scheduled

    # Column C: Tracked (actual) flights.
    $ws.Cells.Item($templateRow, 3).Copy($ws.Cells.Item($row, 3))
    $ws.Cells.Item($row, 3).Value = $r.Tracked

    # Column D: Percent on-time = Tracked / Scheduled.
    $ws.Cells.Item($templateRow, 4).Copy($ws.Cells.Item($row, 4))
    $ws.Cells.Item($row, 4).Formula = "=C$row/B$row"
}

# Matches the author's final selection in the saved workbook.
$ws.Range("G655").Select()
